# Update the "想去人数" (F column) figures on both the "展览" and
# "全部类型" worksheets to reflect the newly scraped counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 170
$ws1.Range("F4").Value  = 171
$ws1.Range("F5").Value  = 4855
$ws1.Range("F6").Value  = 23
$ws1.Range("F7").Value  = 44
$ws1.Range("F8").Value  = 4
$ws1.Range("F10").Value = 487
$ws1.Range("F13").Value = 1360
$ws1.Range("F14").Value = 3273
$ws1.Range("F16").Value = 119
$ws1.Range("F17").Value = 101
$ws1.Range("F18").Value = 70
$ws1.Range("F19").Value = 2511
$ws1.Range("F20").Value = 119
$ws1.Range("F23").Value = 172
$ws1.Range("F24").Value = 31
$ws1.Range("F25").Value = 121
$ws1.Range("F27").Value = 249

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 170
$ws4.Range("F4").Value  = 171
$ws4.Range("F6").Value  = 4855
$ws4.Range("F7").Value  = 23
$ws4.Range("F8").Value  = 44
$ws4.Range("F9").Value  = 4
$ws4.Range("F11").Value = 487
$ws4.Range("F14").Value = 1360
$ws4.Range("F15").Value = 3274
$ws4.Range("F17").Value = 119
$ws4.Range("F18").Value = 101
$ws4.Range("F19").Value = 70
$ws4.Range("F20").Value = 2511
$ws4.Range("F21").Value = 119
$ws4.Range("F24").Value = 172
$ws4.Range("F25").Value = 31
$ws4.Range("F26").Value = 121
$ws4.Range("F28").Value = 249
